$d = $word.ActiveDocument

# 1. Update the letter date: September 19, 2025 -> September 21, 2025
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false,
                         $true, 1, $false, "September 21, 2025", 2)

# 2. Split the single-line mailing address "999 Story Road, San Jose CA 95122"
#    into two paragraphs: "999 Story Road" followed by a new paragraph
#    "San Jose, CA 95122" (same paragraph/run formatting as the original).
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -match "^999 Story Road, San Jose CA 95122\r?$") {
        # Duplicate the paragraph (copies pPr/rPr) right after this one.
        $p.Range.InsertParagraphAfter()
        $newPara = $d.Paragraphs($i + 1)
        $newPara.Range.Text = "San Jose, CA 95122"

        # Trim the original paragraph down to just the street address.
        $p.Range.Find.Execute(", San Jose CA 95122", $true, $false, $false, $false, $false,
                               $true, 1, $false, "", 2)
        break
    }
}

# 3. Remove the blank "No Spacing" paragraph that directly follows
#    "Vietnam Town Condominium Owners Association Board of Directors".
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -match "Board of Directors") {
        $next = $d.Paragraphs($i + 1)
        if (($next.Range.Text -match "^\s*$") -and ($next.Style.NameLocal -eq "No Spacing")) {
            $next.Range.Delete()
        }
        break
    }
}
